# This script applies the latest cryptocurrency market snapshot (price and
# 1h volume change) to the "cryptos" worksheet, including a couple of rows
# whose coins swapped rank position in the source ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that keeps the workbook's default (unformatted) style. We reuse its
# style below to restore formatting on cells whose value we had to coerce to
# text, since numeric-looking strings (e.g. "253.87") would otherwise be
# auto-converted to actual numbers by the Value setter.
$defaultStyle = $ws.Range("B2").Style

function Set-CellText($ws, $cellRef, $value, $forceText, $defaultStyle) {
    $range = $ws.Range($cellRef)
    if ($forceText) {
        # Prepend a text-prefix marker so the numeric-looking value is stored
        # as a string instead of being parsed into a number, then restore the
        # cell's original (default) style/number format.
        $range.Value = "'" + $value
        $range.Style = $defaultStyle
    } else {
        $range.Value = $value
    }
}

$updates = @(
    @{ Cell = "D2"; Value = '37.593.59'; ForceText = $false },
    @{ Cell = "E2"; Value = '  +5.57%  '; ForceText = $false },
    @{ Cell = "D3"; Value = '2.062.26'; ForceText = $false },
    @{ Cell = "E4"; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = "D5"; Value = '253.87'; ForceText = $true },
    @{ Cell = "E5"; Value = '  +3.12%  '; ForceText = $false },
    @{ Cell = "E6"; Value = '  +3.31%  '; ForceText = $false },
    @{ Cell = "D7"; Value = '67.70'; ForceText = $true },
    @{ Cell = "E7"; Value = '  +16.54%  '; ForceText = $false },
    @{ Cell = "E8"; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = "D9"; Value = '0.393'; ForceText = $true },
    @{ Cell = "E9"; Value = '  +8.81%  '; ForceText = $false },
    @{ Cell = "D10"; Value = '59.94'; ForceText = $true },
    @{ Cell = "E10"; Value = '  +2.09%  '; ForceText = $false },
    @{ Cell = "D11"; Value = '0.0774'; ForceText = $true },
    @{ Cell = "E11"; Value = '  +5.05%  '; ForceText = $false },
    @{ Cell = "E12"; Value = '  +0.59%  '; ForceText = $false },
    @{ Cell = "D13"; Value = '0.940'; ForceText = $true },
    @{ Cell = "E13"; Value = '  -2.20%  '; ForceText = $false },
    @{ Cell = "D14"; Value = '23.99'; ForceText = $true },
    @{ Cell = "E14"; Value = '  +29.69%  '; ForceText = $false },
    @{ Cell = "D15"; Value = '15.08'; ForceText = $true },
    @{ Cell = "E15"; Value = '  +3.29%  '; ForceText = $false },
    @{ Cell = "D16"; Value = '2.364.91'; ForceText = $false },
    @{ Cell = "E16"; Value = '  +4.10%  '; ForceText = $false },
    @{ Cell = "D17"; Value = '5.69'; ForceText = $true },
    @{ Cell = "E17"; Value = '  +7.63%  '; ForceText = $false },
    @{ Cell = "D18"; Value = '2.065.23'; ForceText = $false },
    @{ Cell = "E18"; Value = '  +4.27%  '; ForceText = $false },
    @{ Cell = "D19"; Value = '37.522.72'; ForceText = $false },
    @{ Cell = "E19"; Value = '  +5.50%  '; ForceText = $false },
    @{ Cell = "D20"; Value = '73.84'; ForceText = $true },
    @{ Cell = "E20"; Value = '  +3.23%  '; ForceText = $false },
    @{ Cell = "D21"; Value = '0.0₃0882'; ForceText = $false },
    @{ Cell = "E21"; Value = '  +4.03%  '; ForceText = $false },
    @{ Cell = "D22"; Value = '5.52'; ForceText = $true },
    @{ Cell = "E22"; Value = '  +5.41%  '; ForceText = $false },
    @{ Cell = "D23"; Value = '241.01'; ForceText = $true },
    @{ Cell = "E23"; Value = '  +3.53%  '; ForceText = $false },
    @{ Cell = "D24"; Value = '2.70'; ForceText = $true },
    @{ Cell = "E24"; Value = '  +5.05%  '; ForceText = $false },
    @{ Cell = "E25"; Value = '  +0.05%  '; ForceText = $false },
    @{ Cell = "D26"; Value = '2.45'; ForceText = $true },
    @{ Cell = "E26"; Value = '  +7.38%  '; ForceText = $false },
    @{ Cell = "D27"; Value = '10.16'; ForceText = $true },
    @{ Cell = "E27"; Value = '  +10.96%  '; ForceText = $false },
    @{ Cell = "D28"; Value = '162.77'; ForceText = $true },
    @{ Cell = "E28"; Value = '  -1.27%  '; ForceText = $false },
    @{ Cell = "B29"; Value = 'Kaspa'; ForceText = $false },
    @{ Cell = "C29"; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText = $false },
    @{ Cell = "D29"; Value = '0.136'; ForceText = $true },
    @{ Cell = "E29"; Value = '  +43.11%  '; ForceText = $false },
    @{ Cell = "B30"; Value = 'EthereumClassic'; ForceText = $false },
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText = $false },
    @{ Cell = "D30"; Value = '20.21'; ForceText = $true },
    @{ Cell = "E30"; Value = '  +5.06%  '; ForceText = $false },
    @{ Cell = "E31"; Value = '  +3.34%  '; ForceText = $false },
    @{ Cell = "E32"; Value = '  +7.55%  '; ForceText = $false },
    @{ Cell = "D33"; Value = '1.22'; ForceText = $true },
    @{ Cell = "E33"; Value = '  +10.51%  '; ForceText = $false },
    @{ Cell = "D34"; Value = '4.77'; ForceText = $true },
    @{ Cell = "E34"; Value = '  +9.05%  '; ForceText = $false },
    @{ Cell = "D35"; Value = '0.0635'; ForceText = $true },
    @{ Cell = "E35"; Value = '  +6.57%  '; ForceText = $false },
    @{ Cell = "D36"; Value = '2.45'; ForceText = $true },
    @{ Cell = "E36"; Value = '  +0.02%  '; ForceText = $false },
    @{ Cell = "D37"; Value = '6.34'; ForceText = $true },
    @{ Cell = "E37"; Value = '  +16.29%  '; ForceText = $false },
    @{ Cell = "E38"; Value = '  +0.14%  '; ForceText = $false },
    @{ Cell = "E39"; Value = '  +3.05%  '; ForceText = $false },
    @{ Cell = "D40"; Value = '3.15'; ForceText = $true },
    @{ Cell = "E40"; Value = '  +38.50%  '; ForceText = $false },
    @{ Cell = "D41"; Value = '0.104'; ForceText = $true },
    @{ Cell = "E41"; Value = '  +14.75%  '; ForceText = $false },
    @{ Cell = "B42"; Value = 'InjectiveProtocol'; ForceText = $false },
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false },
    @{ Cell = "D42"; Value = '18.43'; ForceText = $true },
    @{ Cell = "E42"; Value = '  +13.90%  '; ForceText = $false },
    @{ Cell = "B43"; Value = 'TrustWalletToken'; ForceText = $false },
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false },
    @{ Cell = "D43"; Value = '1.28'; ForceText = $true },
    @{ Cell = "E43"; Value = '  +3.48%  '; ForceText = $false },
    @{ Cell = "B44"; Value = 'HuobiToken'; ForceText = $false },
    @{ Cell = "C44"; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; ForceText = $false },
    @{ Cell = "D44"; Value = '3.06'; ForceText = $true },
    @{ Cell = "E44"; Value = '  +6.45%  '; ForceText = $false },
    @{ Cell = "B45"; Value = 'VeChain'; ForceText = $false },
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false },
    @{ Cell = "D45"; Value = '0.0221'; ForceText = $true },
    @{ Cell = "E45"; Value = '  +3.87%  '; ForceText = $false },
    @{ Cell = "B46"; Value = 'ARBITRUM'; ForceText = $false },
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false },
    @{ Cell = "D46"; Value = '1.16'; ForceText = $true },
    @{ Cell = "E46"; Value = '  +5.74%  '; ForceText = $false },
    @{ Cell = "D47"; Value = '97.73'; ForceText = $true },
    @{ Cell = "E47"; Value = '  +4.47%  '; ForceText = $false },
    @{ Cell = "D48"; Value = '8.01'; ForceText = $true },
    @{ Cell = "E48"; Value = '  +2.77%  '; ForceText = $false },
    @{ Cell = "D49"; Value = '1.422.88'; ForceText = $false },
    @{ Cell = "E49"; Value = '  +3.57%  '; ForceText = $false },
    @{ Cell = "D50"; Value = '2.96'; ForceText = $true },
    @{ Cell = "E50"; Value = '  +2.02%  '; ForceText = $false },
    @{ Cell = "D51"; Value = '3.78'; ForceText = $true },
    @{ Cell = "E51"; Value = '  +9.65%  '; ForceText = $false }
)

foreach ($u in $updates) {
    Set-CellText $ws $u.Cell $u.Value $u.ForceText $defaultStyle
}
